$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for the specified rows to match re-pulled data
$ws.Range("F7").Value = -2
$ws.Range("F8").Value = 0
$ws.Range("F11").Value = -1
$ws.Range("F12").Value = 1
$ws.Range("F15").Value = -1
$ws.Range("F21").Value = 2
$ws.Range("F24").Value = 3
$ws.Range("F32").Value = -2
$ws.Range("F34").Value = 1
$ws.Range("F36").Value = 1
$ws.Range("F38").Value = -2
$ws.Range("F41").Value = -5
$ws.Range("F42").Value = -2
$ws.Range("F48").Value = -2
$ws.Range("F54").Value = -2
$ws.Range("F56").Value = 5
$ws.Range("F59").Value = -3
$ws.Range("F62").Value = -2
$ws.Range("F64").Value = -1
